$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells that look like plain numbers stay formatted as text,
# matching the original inline-string representation (preserves exact digits).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values per the diff
$ws.Range("D2").Value = '42.284.12'
$ws.Range("E2").Value = '  -0.76%  '

$ws.Range("D3").Value = '2.176.73'
$ws.Range("E3").Value = '  -1.86%  '

$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").Value = '251.90'
$ws.Range("E5").Value = '  +4.71%  '

$ws.Range("D6").Value = '0.610'
$ws.Range("E6").Value = '  -1.58%  '

$ws.Range("D7").Value = '72.84'
$ws.Range("E7").Value = '  -2.55%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '0.579'
$ws.Range("E9").Value = '  -3.24%  '

$ws.Range("D10").Value = '39.90'
$ws.Range("E10").Value = '  -3.01%  '

$ws.Range("E11").Value = '  -2.41%  '

$ws.Range("E12").Value = '  -0.58%  '

$ws.Range("E13").Value = '  -2.47%  '

$ws.Range("D14").Value = '2.502.45'
$ws.Range("E14").Value = '  -1.98%  '

$ws.Range("D15").Value = '14.14'
$ws.Range("E15").Value = '  -3.62%  '

$ws.Range("D16").Value = '2.173.83'
$ws.Range("E16").Value = '  -2.28%  '

$ws.Range("E17").Value = '  -3.97%  '

$ws.Range("D18").Value = '42.191.94'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").Value = '0.0000102'
$ws.Range("E19").Value = '  -3.10%  '

$ws.Range("D20").Value = '70.49'
$ws.Range("E20").Value = '  -0.38%  '

$ws.Range("E21").Value = '  -1.81%  '

$ws.Range("D22").Value = '226.40'
$ws.Range("E22").Value = '  -1.29%  '

$ws.Range("D23").Value = '9.27'
$ws.Range("E23").Value = '  -5.79%  '

$ws.Range("E24").Value = '  -2.28%  '

$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("D26").Value = '10.40'
$ws.Range("E26").Value = '  -4.58%  '

$ws.Range("E27").Value = '  -0.30%  '

$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.16'
$ws.Range("E28").Value = '  -2.09%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.13'
$ws.Range("E29").Value = '  -2.33%  '

$ws.Range("D30").Value = '170.13'
$ws.Range("E30").Value = '  -1.57%  '

$ws.Range("D31").Value = '36.33'
$ws.Range("E31").Value = '  -0.17%  '

$ws.Range("D32").Value = '19.92'
$ws.Range("E32").Value = '  -1.70%  '

$ws.Range("D33").Value = '0.0809'
$ws.Range("E33").Value = '  +1.70%  '

$ws.Range("D34").Value = '5.07'
$ws.Range("E34").Value = '  -4.18%  '

$ws.Range("D35").Value = '0.119'
$ws.Range("E35").Value = '  -1.74%  '

$ws.Range("E36").Value = '  -2.51%  '

$ws.Range("D37").Value = '0.0334'
$ws.Range("E37").Value = '  +4.06%  '

$ws.Range("E38").Value = '  -5.30%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '2.03'
$ws.Range("E39").Value = '  -4.62%  '

$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = '11.62'
$ws.Range("E40").Value = '  -6.17%  '

$ws.Range("E41").Value = '  -1.53%  '

$ws.Range("D42").Value = '58.85'
$ws.Range("E42").Value = '  -2.49%  '

$ws.Range("D43").Value = '5.12'
$ws.Range("E43").Value = '  -6.88%  '

$ws.Range("D44").Value = '100.81'
$ws.Range("E44").Value = '  +1.62%  '

$ws.Range("E45").Value = '  +6.90%  '

$ws.Range("D46").Value = '0.0966'
$ws.Range("E46").Value = '  -2.42%  '

$ws.Range("D47").Value = '8.14'
$ws.Range("E47").Value = '  -5.12%  '

$ws.Range("E48").Value = '  +3.84%  '

$ws.Range("D49").Value = '1.08'
$ws.Range("E49").Value = '  -2.80%  '

$ws.Range("E50").Value = '  -1.62%  '

$ws.Range("D51").Value = '2.65'
$ws.Range("E51").Value = '  +0.31%  '
